$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from the default "Sheet1" to "AddCustomerTest"
$ws.Name = "AddCustomerTest"

# Header row + first data row for the first three columns (firstname/lastname/postcode)
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

$ws.Range("A2").Value = "thamil"
$ws.Range("B2").Value = "alagan"
$ws.Range("C2").Value = "b27 a93"

# alerttext column added afterwards (matches shared-string insertion order in the workbook)
$ws.Range("D1").Value = "alerttext"
$ws.Range("D2").Value = "Customer added successfully"

# Column D was widened (best-fit) to accommodate the longer alert text
$ws.Columns("D:D").AutoFit()
